# "Product Backlog" sheet: drop the old construction-app backlog rows,
# keep only the two Myst game-collection stories, and rename the first
# story from "Add games..." to "Select games...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove backlog rows 5-11 (setup build stages ... share calendar), shifting
# everything below them up so the sheet ends with just rows 1-4.
$ws.Range("A5:D11").EntireRow.Delete()

# Row 2: rename the "Add games that I own/enjoy" story.
$ws.Range("C2").Value = "Select games that I own/enjoy"

# Row 4 previously held the "setup build stages" story; only the leftover
# "User" tag in column B remains, the rest of the row is cleared.
$ws.Range("A4:D4").ClearContents()
$ws.Range("B4").Value = "User"

# Match the saved selection/view state.
$ws.Range("C5").Select()
